$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 5 ("nan" polite_expressions, DIS/WRI review of the
# WordNet/embedding comment) keeps its D:I content, but its
# polite_expressions (C) value is cleared.
$ws.Range("C5").Value = ""

# A new row 6 is appended. Its polite_expressions value ("nan") is what
# used to live in C5, and the rest of the row holds the new DeePa review.
$ws.Range("A6").Value = "parisk"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "nan"
$ws.Range("D6").Value = "DIS"
$ws.Range("E6").Value = "THE"
$ws.Range("F6").Value = "42be9703-0e9b-4ce8-962d-60bf1f233ce8"
$ws.Range("G6").Value = "SJCPLLpaW_annotated.xlsx"
$ws.Range("H6").Value = "The results show that DeePa achieves speedups compared to PyTorch and TensorFlow with all of the tested minibatch sizes."
$ws.Range("I6").Value = "Needs Revision"
